$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# ---------------------------------------------------------------------------
# Row 26 : new "ATOR(ES)" / "CASOS DE USO" header cells (bold, centered)
# Style creation order (B26 -> s8 wrap, C26 -> s9 no-wrap) is decoupled from
# the value-assignment order below (C26 first) that drives shared-string
# insertion order.
# ---------------------------------------------------------------------------
$ws.Cells.Item(26, 2).HorizontalAlignment = -4108
$ws.Cells.Item(26, 2).VerticalAlignment = -4108
$ws.Cells.Item(26, 2).WrapText = $true
$ws.Cells.Item(26, 2).Font.Bold = $true

$ws.Cells.Item(26, 3).HorizontalAlignment = -4108
$ws.Cells.Item(26, 3).VerticalAlignment = -4108
$ws.Cells.Item(26, 3).Font.Bold = $true

$ws.Cells.Item(26, 3).Value = "ATOR(ES)"
$ws.Cells.Item(26, 2).Value = "CASOS DE USO"

# ---------------------------------------------------------------------------
# Column A : UC01..UC10, centered horizontally, top vertically
# ---------------------------------------------------------------------------
$ws.Cells.Item(27, 1).HorizontalAlignment = -4108
$ws.Cells.Item(27, 1).VerticalAlignment = -4160
$ws.Cells.Item(27, 1).Value = "UC01"

$ws.Cells.Item(28, 1).HorizontalAlignment = -4108
$ws.Cells.Item(28, 1).VerticalAlignment = -4160
$ws.Cells.Item(28, 1).Value = "UC02"

$ws.Cells.Item(29, 1).HorizontalAlignment = -4108
$ws.Cells.Item(29, 1).VerticalAlignment = -4160
$ws.Cells.Item(29, 1).Value = "UC03"

$ws.Cells.Item(30, 1).HorizontalAlignment = -4108
$ws.Cells.Item(30, 1).VerticalAlignment = -4160
$ws.Cells.Item(30, 1).Value = "UC04"

$ws.Cells.Item(31, 1).HorizontalAlignment = -4108
$ws.Cells.Item(31, 1).VerticalAlignment = -4160
$ws.Cells.Item(31, 1).Value = "UC05"

$ws.Cells.Item(32, 1).HorizontalAlignment = -4108
$ws.Cells.Item(32, 1).VerticalAlignment = -4160
$ws.Cells.Item(32, 1).Value = "UC06"

$ws.Cells.Item(33, 1).HorizontalAlignment = -4108
$ws.Cells.Item(33, 1).VerticalAlignment = -4160
$ws.Cells.Item(33, 1).Value = "UC07"

$ws.Cells.Item(34, 1).HorizontalAlignment = -4108
$ws.Cells.Item(34, 1).VerticalAlignment = -4160
$ws.Cells.Item(34, 1).Value = "UC08"

$ws.Cells.Item(35, 1).HorizontalAlignment = -4108
$ws.Cells.Item(35, 1).VerticalAlignment = -4160
$ws.Cells.Item(35, 1).Value = "UC09"

$ws.Cells.Item(36, 1).HorizontalAlignment = -4108
$ws.Cells.Item(36, 1).VerticalAlignment = -4160
$ws.Cells.Item(36, 1).Value = "UC10"

# ---------------------------------------------------------------------------
# Column B : use case names (set in the author's original typing order)
# ---------------------------------------------------------------------------
$ws.Cells.Item(30, 2).WrapText = $true
$ws.Cells.Item(30, 2).Value = "Autenticar Veículo"

$ws.Cells.Item(28, 2).Value = "Registrar Saída de Veículo"

$ws.Cells.Item(27, 2).Value = "Registrar Entrada de Veículo"

$ws.Cells.Item(29, 2).VerticalAlignment = -4160
$ws.Cells.Item(29, 2).WrapText = $true
$ws.Cells.Item(29, 2).Value = "Registrar Veículo"

$ws.Cells.Item(31, 2).VerticalAlignment = -4160
$ws.Cells.Item(31, 2).WrapText = $true
$ws.Cells.Item(31, 2).Value = "Registrar Funcionário"

$ws.Cells.Item(32, 2).VerticalAlignment = -4160
$ws.Cells.Item(32, 2).WrapText = $true
$ws.Cells.Item(32, 2).Value = "Excluir Veículo"

$ws.Cells.Item(33, 2).Value = "Manter Vagas"

$ws.Cells.Item(34, 2).WrapText = $true
$ws.Cells.Item(34, 2).Value = "Efetuar Consultas"

$ws.Cells.Item(35, 2).VerticalAlignment = -4160
$ws.Cells.Item(35, 2).Value = "Autenticar Usuário"

$ws.Cells.Item(36, 2).VerticalAlignment = -4160
$ws.Cells.Item(36, 2).Value = "Gerenciar Permissões"

# ---------------------------------------------------------------------------
# Column C : actor(s) (set in the author's original typing order)
# ---------------------------------------------------------------------------
$ws.Cells.Item(27, 3).Value = "Cancela Eletrônica"
$ws.Cells.Item(28, 3).Value = "Cancela Eletrônica"
$ws.Cells.Item(30, 3).Value = "Cancela Eletrônica"

$ws.Cells.Item(32, 3).Value = "Funcionário, Gestor"

$ws.Cells.Item(31, 3).Value = "Setor Pessoal"

$ws.Cells.Item(33, 3).Value = "Gestor"
$ws.Cells.Item(34, 3).Value = "Gestor"
$ws.Cells.Item(36, 3).Value = "Gestor"

$ws.Cells.Item(35, 3).Value = "Gestor, Funcionário, Setor Administrativo, Setor Pessoal"

$ws.Cells.Item(29, 3).Value = "Funcionário, Gestor, Setor Administrativo"

# ---------------------------------------------------------------------------
# Rows 37-40 : trailing formatted-but-empty cells in column B
# ---------------------------------------------------------------------------
$ws.Cells.Item(37, 2).VerticalAlignment = -4160
$ws.Cells.Item(37, 2).WrapText = $true

$ws.Cells.Item(38, 2).VerticalAlignment = -4160

$ws.Cells.Item(39, 2).WrapText = $true

$ws.Cells.Item(40, 2).VerticalAlignment = -4160
$ws.Cells.Item(40, 2).WrapText = $true

# ---------------------------------------------------------------------------
# View state : zoom + active cell selection
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 125
$ws.Range("C36").Select()
